# 汽车.xlsx update:
#  1. The data rows are grouped in 4-row year blocks (A/B/C/D quarters),
#     starting at row 2 (row 1 is the header). Within every block the "B"
#     quarter row and the "C" quarter row have had their A:E contents
#     swapped (the "A" and "D" rows are untouched).
#  2. Columns F (汽车产销率) and G (汽车销售量) - including their header
#     cells in row 1 - are removed entirely, shrinking the used range from
#     A1:G81 to A1:E81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($block = 0; $block -lt 20; $block++) {
    $base = 2 + $block * 4   # first row ("A" quarter) of this year's block
    $rowB = $base + 1        # "B" quarter row
    $rowC = $base + 2        # "C" quarter row

    $valuesB = $ws.Range("A$($rowB):E$($rowB)").Value2
    $valuesC = $ws.Range("A$($rowC):E$($rowC)").Value2

    $ws.Range("A$($rowB):E$($rowB)").Value = $valuesC
    $ws.Range("A$($rowC):E$($rowC)").Value = $valuesB
}

# Delete the now-unwanted F:G columns (产销率 / 销售量 duplicates).
$ws.Columns("F:G").Delete()
